# Applies the "Wentworth grade scale" paragraph block after the existing
# image paragraph, and marks the image run as NoProof (<w:noProof/>).

$d = $word.ActiveDocument

# --- 1. Append 11 new empty paragraphs after the image paragraph ---
$imgPara = $d.Paragraphs(1)
$tail = $imgPara.Range
for ($i = 0; $i -lt 11; $i++) {
    $tail.InsertParagraphAfter() | Out-Null
}

# --- 2. Mark the drawing's run as NoProof (adds <w:rPr><w:noProof/></w:rPr>) ---
#        Done after the new paragraphs exist so the property doesn't bleed
#        into them.
$d.Paragraphs(1).Range.NoProofing = 1

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$rPrGeorgia = '<w:rPr><w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/><w:color w:val="1F1F1F"/></w:rPr>'
$pPrGeorgia = "<w:pPr>$rPrGeorgia</w:pPr>"

# Paragraphs 2 & 3: fully empty paragraphs (self-closing <w:p/>, no run at all)
$pEmpty = "<w:p $wNs/>"
$d.Paragraphs(2).Range.InsertXML($pEmpty) | Out-Null
$d.Paragraphs(3).Range.InsertXML($pEmpty) | Out-Null

# Paragraph 4: "Wentworth grade scale "
$p4 = "<w:p $wNs>$pPrGeorgia<w:r>$rPrGeorgia<w:t xml:space=`"preserve`">Wentworth grade scale </w:t></w:r></w:p>"
$d.Paragraphs(4).Range.InsertXML($p4) | Out-Null

# Paragraph 5: empty, but pPr carries the Georgia/1F1F1F run mark
$p5 = "<w:p $wNs>$pPrGeorgia</w:p>"
$d.Paragraphs(5).Range.InsertXML($p5) | Out-Null

# Paragraph 6: gravel-sized particles ...
$p6 = "<w:p $wNs>" +
      "$pPrGeorgia" +
      "<w:r>$rPrGeorgia<w:t xml:space=`"preserve`">gravel-sized particles have a nominal diameter of </w:t></w:r>" +
      '<w:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:color w:val="1F1F1F"/></w:rPr><w:t>&#x2A7E;</w:t></w:r>' +
      "<w:r>$rPrGeorgia<w:t xml:space=`"preserve`">2.0 mm; </w:t></w:r>" +
      "</w:p>"
$d.Paragraphs(6).Range.InsertXML($p6) | Out-Null

# Paragraph 7: sand-sized particles ...
$p7 = "<w:p $wNs>" +
      "$pPrGeorgia" +
      "<w:r>$rPrGeorgia<w:t xml:space=`"preserve`">sand-sized particles have nominal diameters from &lt;2.0 mm to </w:t></w:r>" +
      '<w:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:color w:val="1F1F1F"/></w:rPr><w:t>&#x2A7E;</w:t></w:r>' +
      "<w:r>$rPrGeorgia<w:t>62.5 </w:t></w:r>" +
      '<w:proofErr w:type="spellStart"/>' +
      "<w:r>$rPrGeorgia<w:t>&#x3BC;m</w:t></w:r>" +
      '<w:proofErr w:type="spellEnd"/>' +
      "<w:r>$rPrGeorgia<w:t>;</w:t></w:r>" +
      "</w:p>"
$d.Paragraphs(7).Range.InsertXML($p7) | Out-Null

# Paragraph 8: silt-sized particles ...
$p8 = "<w:p $wNs>" +
      "$pPrGeorgia" +
      "<w:r>$rPrGeorgia<w:t xml:space=`"preserve`">silt-sized particles have nominal diameters from &lt;62.5 to </w:t></w:r>" +
      '<w:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:color w:val="1F1F1F"/></w:rPr><w:t>&#x2A7E;</w:t></w:r>' +
      "<w:r>$rPrGeorgia<w:t>4.0 </w:t></w:r>" +
      '<w:proofErr w:type="spellStart"/>' +
      "<w:r>$rPrGeorgia<w:t>&#x3BC;m</w:t></w:r>" +
      '<w:proofErr w:type="spellEnd"/>' +
      "<w:r>$rPrGeorgia<w:t xml:space=`"preserve`">; and </w:t></w:r>" +
      "</w:p>"
$d.Paragraphs(8).Range.InsertXML($p8) | Out-Null

# Paragraph 9: clay is <4.0 μm
$p9 = "<w:p $wNs>" +
      "$pPrGeorgia" +
      "<w:r>$rPrGeorgia<w:t>clay is &lt;4.0 </w:t></w:r>" +
      '<w:proofErr w:type="spellStart"/>' +
      "<w:r>$rPrGeorgia<w:t>&#x3BC;m</w:t></w:r>" +
      '<w:proofErr w:type="spellEnd"/>' +
      "</w:p>"
$d.Paragraphs(9).Range.InsertXML($p9) | Out-Null

# Paragraph 10 & 11: empty, but pPr carries the Georgia/1F1F1F run mark
$p10 = "<w:p $wNs>$pPrGeorgia</w:p>"
$d.Paragraphs(10).Range.InsertXML($p10) | Out-Null
$d.Paragraphs(11).Range.InsertXML($p10) | Out-Null

# Paragraph 12: "Shepard (1954)" (no pPr, just the run rPr)
$p12 = "<w:p $wNs><w:r>$rPrGeorgia<w:t>Shepard (1954)</w:t></w:r></w:p>"
$d.Paragraphs(12).Range.InsertXML($p12) | Out-Null

Write-Output "done"
